$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H ("input shape")
$ws.Range("H1").Value = "input shape"

# --- Row 2 ---
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 0.01648409292101860046
$ws.Range("D2").Value = 0.0941374972462653975
$ws.Range("E2").Value = 0.12839040160179099326
$ws.Range("F2").Value = 0.00808862876147030917
$ws.Range("G2").Value = 0.58377587795257501746
$ws.Range("C2:G2").NumberFormat = "0.000"
$ws.Range("H2").Value = 256
$ws.Range("H2").NumberFormat = "General"

# --- Row 3 ---
$ws.Range("B3").Value = 9
$ws.Range("C3:G3").ClearContents()
$ws.Range("H3").Value = 300
$ws.Range("H3").NumberFormat = "General"

# --- Row 4 --- (A4 stays 3; clear the rest of the old row)
$ws.Range("B4:G4").ClearContents()

# Column widths for C:G (matches bestFit width Excel computed for the numbers)
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 11.6640625

# Selection ends on D5 after the edit
$ws.Range("D5").Select()
